$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "final row" date number format (currently on A6) before we move it.
$lastRowFormat = $ws.Range("A6").NumberFormat

# Previous last row (row 6) reverts to the normal date-time format used by
# all other data rows (same format as the rest of column A).
$ws.Range("A6").NumberFormat = $ws.Range("A2").NumberFormat

# Append the new day's row (row 7) with the "final row" date format.
$ws.Range("A7").Value = 45956
$ws.Range("A7").NumberFormat = $lastRowFormat
$ws.Range("B7").Value = 12
$ws.Range("C7").Value = 16
$ws.Range("D7").Value = 14
